$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet from "Plan1" to "DOCUMENTOS"
$ws.Name = "DOCUMENTOS"

# Reset selection to A1 (removes the stored selection at E16)
$ws.Range("A1").Select()
